# Commit: "Adicionando capacidade ao projeto para trabalhar com relatorios allure"
#
# Update the scenario results table on the active sheet:
#   - C2 (Status for CT 01):   "Failed" -> "Passed"
#   - H2 (vOutData for CT 01): "16/06/2020" -> "16/04/2021"
#   - B3 (RunTest for CT 02):  "No" -> "Yes"
#   - H3 (vOutData for CT 02): "26/05/2020" -> "16/04/2021"
# and move the active cell selection to B4 (as last saved by Excel).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "Passed"
$ws.Range("H2").Value = "16/04/2021"
$ws.Range("B3").Value = "Yes"
$ws.Range("H3").Value = "16/04/2021"

[void]$ws.Range("B4").Select()
